$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.895.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.98%  "
$ws.Range("D3").Value = "'2.635.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'528.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.91%  "
$ws.Range("D6").Value = "'155.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("E10").Value = "  +5.15%  "
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'3.097.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").Value = "'60.891.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.98%  "
$ws.Range("D15").Value = "'21.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "'0.0000143"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "'2.644.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").Value = "'4.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'353.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("D21").Value = "'6.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.67%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").Value = "'61.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("E24").Value = "  +2.50%  "
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").Value = "'0.977"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("D27").Value = "'0.0₃0863"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.74%  "
$ws.Range("D28").Value = "'7.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'6.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.87%  "
$ws.Range("D31").Value = "'19.45"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  +3.99%  "
$ws.Range("D33").Value = "'150.15"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.96%  "
$ws.Range("E34").Value = "  +3.92%  "
$ws.Range("E35").Value = "  +1.67%  "
$ws.Range("D36").Value = "'0.926"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.72%  "
$ws.Range("D37").Value = "'0.894"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.51%  "
$ws.Range("E38").Value = "  +1.49%  "
$ws.Range("D39").Value = "'305.65"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.72%  "
$ws.Range("D40").Value = "'3.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.69%  "
$ws.Range("D41").Value = "'36.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("E42").Value = "  +4.08%  "
$ws.Range("E43").Value = "  +1.73%  "
$ws.Range("E44").Value = "  +1.93%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").Value = "'19.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  +2.94%  "
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").Value = "'19.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.20%  "
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").Value = "'1.982.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.21%  "
